# Update 30-Mei-2021, end of day update.
#
# The daily petty-cash ledger on "Sheet1" is reset for the new day: the
# individual transaction rows (date / description / debit / credit, columns
# A-D) for the previous period are cleared out, and the running "Saldo"
# balance in column E is reset to the new opening value. The existing
# shared running-balance formulas in column E are left untouched, so once
# the debit/credit cells they reference are cleared the whole column simply
# carries the new opening balance straight down through row 114.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New opening ("SALDO AWAL") balance for the day, was 1,124,025.
$ws.Range("E2").Value = 473125

# Row 3 still records its "Wages Expense" line, but it now belongs to a
# later date and no longer carries a debit amount.
$ws.Range("A3").Value = 44347
$ws.Range("D3").Clear()

# Rows 4-42 held the rest of the prior period's dated transactions
# (date/description/debit/credit). Clear all of that out; only the
# column-E running-balance formulas remain.
$ws.Range("A4:D42").Clear()

# Restore the frozen-pane view to the top of the ledger (it had scrolled
# down to row 37 with D58 selected).
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
[void]$ws.Range("C4").Select()
